$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ALC_data = @(
  @("H18", 2214.889),
  @("I18", 1241.75),
  @("K18", 1241.75),
  @("M18", -957.75),
  @("H28", 1439.8788),
  @("I28", 1458.3704),
  @("K28", 1458.3704),
  @("M28", -973.3704),
  @("H43", 90000),
  @("I43", 100000),
  @("J43", 80000),
  @("K43", 100000),
  @("L43", 80000),
  @("M43", -99931),
  @("N43", -80138),
  @("H53", 216.3077),
  @("I53", 198.375),
  @("J53", 245),
  @("K53", 198.375),
  @("L53", 245),
  @("M53", 438.625),
  @("N53", -1519),
  @("H105", 59999.5),
  @("J105", 59999.5),
  @("L105", 59999.5),
  @("N105", -66987.5),
  @("H113", 6927.684),
  @("J113", 18185.334),
  @("L113", 18185.334),
  @("N113", -24693.334),
  @("H132", 1698.2646),
  @("I132", 1684.1746),
  @("J132", 1875.8),
  @("K132", 5052.5238),
  @("L132", 5627.4),
  @("M132", -2522.5238),
  @("N132", -10687.4),
  @("H138", 8336422.5),
  @("I138", 1604.2307),
  @("J138", 10641798),
  @("K138", 4812.6921),
  @("L138", 31925394),
  @("M138", 327.3078999999998),
  @("N138", -31935674)
)
foreach ($row in $ALC_data) {
  $ws.Range($row[0]).Value = $row[1]
}

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ARM_data = @(
  @("H45", 5216.7036),
  @("I45", 6195.476),
  @("K45", 6195.476),
  @("M45", -5818.476),
  @("H61", 255815.22),
  @("I61", 4499.6665),
  @("K61", 4499.6665),
  @("M61", -4287.6665),
  @("H74", 5414.8774),
  @("I74", 1197.8049),
  @("K74", 1197.8049),
  @("M74", -323.8049000000001),
  @("H77", 5414.8774),
  @("I77", 1197.8049),
  @("K77", 5989.0245),
  @("M77", -1621.0245),
  @("H92", 29000),
  @("J92", 29000),
  @("L92", 29000),
  @("N92", -33992),
  @("H105", 78525.3),
  @("J105", 78694.78),
  @("L105", 78694.78),
  @("N105", -85682.78),
  @("H136", 255815.22),
  @("I136", 4499.6665),
  @("K136", 13498.9995),
  @("M136", -10948.9995)
)
foreach ($row in $ARM_data) {
  $ws.Range($row[0]).Value = $row[1]
}

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$BSM_data = @(
  @("H44", 24799.8),
  @("I44", 16999.5),
  @("J44", 30000),
  @("K44", 16999.5),
  @("L44", 30000),
  @("M44", -16502.5),
  @("N44", -30994),
  @("H134", 1825.5938),
  @("I134", 1637.6333),
  @("J134", 4645),
  @("K134", 4912.8999),
  @("L134", 13935),
  @("M134", -2377.8999),
  @("N134", -19005)
)
foreach ($row in $BSM_data) {
  $ws.Range($row[0]).Value = $row[1]
}

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$CRP_data = @(
  @("H7", 1212.7693),
  @("I7", 1631.4286),
  @("K7", 1631.4286),
  @("M7", -1518.4286),
  @("H32", 24663),
  @("I32", 19595.8),
  @("J32", 49999),
  @("K32", 19595.8),
  @("L32", 49999),
  @("M32", -19279.8),
  @("N32", -50631),
  @("H94", 2297.5),
  @("J94", 2594.2856),
  @("L94", 2594.2856),
  @("N94", -3496.2856),
  @("H99", 3011.6875),
  @("I99", 2722.6365),
  @("K99", 2722.6365),
  @("M99", -1224.6365),
  @("H102", 34990),
  @("J102", 34990),
  @("L102", 34990),
  @("N102", -39858),
  @("H107", 979.1111),
  @("I107", 581.6667),
  @("K107", 581.6667),
  @("M107", 1338.3333),
  @("H124", 48997.5),
  @("J124", 48997.5),
  @("L124", 48997.5),
  @("N124", -53907.5),
  @("H126", 3011.6875),
  @("I126", 2722.6365),
  @("K126", 8167.9095),
  @("M126", -5697.9095),
  @("H141", 241179.6),
  @("J141", 241179.6),
  @("L141", 241179.6),
  @("N141", -251539.6)
)
foreach ($row in $CRP_data) {
  $ws.Range($row[0]).Value = $row[1]
}

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$CUL_data = @(
  @("H33", 135.25),
  @("I33", 79.92856999999999),
  @("K33", 479.57142),
  @("M33", -196.57142),
  @("H44", 2119.4443),
  @("I44", 715),
  @("K44", 2145),
  @("M44", -1747),
  @("H124", 45003),
  @("I124", 30),
  @("J124", 50000),
  @("K124", 90),
  @("L124", 150000),
  @("M124", 4820),
  @("N124", -159820),
  @("H126", 7500),
  @("I126", 7500),
  @("K126", 22500),
  @("M126", -17560),
  @("H129", 11001737),
  @("I129", 19800726),
  @("K129", 59402178),
  @("M129", -59397178)
)
foreach ($row in $CUL_data) {
  $ws.Range($row[0]).Value = $row[1]
}

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$GSM_data = @(
  @("H70", 24499),
  @("I70", 6998),
  @("J70", 42000),
  @("K70", 6998),
  @("L70", 42000),
  @("M70", -6728),
  @("N70", -42540),
  @("H73", 24499),
  @("I73", 6998),
  @("J73", 42000),
  @("K73", 6998),
  @("L73", 42000),
  @("M73", -6062),
  @("N73", -43872),
  @("H92", 13375.167),
  @("I92", 10000),
  @("K92", 10000),
  @("M92", -8128),
  @("H132", 4676.387),
  @("I132", 2109.4075),
  @("J132", 22003.5),
  @("K132", 6328.2225),
  @("L132", 66010.5),
  @("M132", -3798.2225),
  @("N132", -71070.5)
)
foreach ($row in $GSM_data) {
  $ws.Range($row[0]).Value = $row[1]
}

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$LTW_data = @(
  @("H46", 410.83334),
  @("I46", 407.27274),
  @("K46", 407.27274),
  @("M46", -219.27274),
  @("H93", 1716.2667),
  @("I93", 1746.32),
  @("K93", 1746.32),
  @("M93", -498.3199999999999),
  @("H136", 3168.4119),
  @("I136", 2804.5),
  @("J136", 4866.6665),
  @("K136", 8413.5),
  @("L136", 14599.9995),
  @("M136", -5863.5),
  @("N136", -19699.9995)
)
foreach ($row in $LTW_data) {
  $ws.Range($row[0]).Value = $row[1]
}

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$WVR_data = @(
  @("H22", 5003820),
  @("I22", 20000000),
  @("K22", 20000000),
  @("M22", -19999707),
  @("H52", 20315.5),
  @("J52", 0),
  @("L52", 0),
  @("H95", 36602.5),
  @("J95", 36602.5),
  @("L95", 36602.5),
  @("N95", -42094.5),
  @("H130", 29999),
  @("J130", 29999),
  @("L130", 29999),
  @("N130", -40039),
  @("H136", 2835.8096),
  @("I136", 1811.1428),
  @("K136", 5433.428400000001),
  @("M136", -2883.428400000001)
)
foreach ($row in $WVR_data) {
  $ws.Range($row[0]).Value = $row[1]
}
$ws.Range("N52").ClearContents()

Write-Output "Applied all Leve profit updates."